$wb = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item("20111226")
$ws4.Activate()
Write-Output "done"
